$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, copying the header style from E1 (bold, bordered, centered)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill new boolean (FALSE) columns for rows 2-8
$ws.Range("F2:H8").Value = $false
